$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing cells used purely as formatting templates (copy format only, values untouched).
# Row 166 col C already has style s="6" (Arial 11 / vertical-center xf).
# Row 45  col C already has style s="5" (Arial 11 / vertical-center xf).
$styleSix = $ws.Cells.Item(166, 3)
$styleFive = $ws.Cells.Item(45, 3)

# ---- Row 168 ----
$c168 = $ws.Cells.Item(168, 3)
$prefix168 = "Southern California Academy of Science: "
$suffix168 = "Risk Characterization of Microplastics in San Francisco Bay, California"
$c168.Value2 = $prefix168 + $suffix168
$pr168 = $c168.Characters(1, $prefix168.Length)
$pr168.Font.Name = "Arial"
$pr168.Font.Size = 11
$pr168.Font.Color = 0
$pr168.Font.Italic = $false
$it168 = $c168.Characters($prefix168.Length + 1, $suffix168.Length)
$it168.Font.Name = "Arial"
$it168.Font.Size = 11
$it168.Font.Color = 0
$it168.Font.Italic = $true
$styleSix.Copy()
$c168.PasteSpecial(-4122)

$ws.Cells.Item(168, 1).Value2 = 2022
$ws.Cells.Item(168, 2).Value2 = "May"
$ws.Cells.Item(168, 4).Value2 = "Platform"
$ws.Cells.Item(168, 5).Value2 = "in person"
$ws.Cells.Item(168, 6).Value2 = "conference"

# ---- Row 169 ----
$c169 = $ws.Cells.Item(169, 3)
$prefix169 = "U.S. Government Interest Group on Nanoplastics: "
$suffix169 = "California's Actions on Microplastics"
$c169.Value2 = $prefix169 + $suffix169
$it169 = $c169.Characters($prefix169.Length + 1, $suffix169.Length)
$it169.Font.Name = "Arial"
$it169.Font.Size = 11
$it169.Font.Color = 0
$it169.Font.Italic = $true
$styleFive.Copy()
$c169.PasteSpecial(-4122)

$ws.Cells.Item(169, 1).Value2 = 2022
$ws.Cells.Item(169, 2).Value2 = "May"
$ws.Cells.Item(169, 4).Value2 = "Platform"
$ws.Cells.Item(169, 5).Value2 = "virtual"
$ws.Cells.Item(169, 6).Value2 = "meeting"

# ---- Row 170 ----
$c170 = $ws.Cells.Item(170, 3)
$prefix170 = "AARP: "
$suffix170 = "Plastics and Your Health"
$c170.Value2 = $prefix170 + $suffix170
$pr170 = $c170.Characters(1, $prefix170.Length)
$pr170.Font.Name = "Arial"
$pr170.Font.Size = 11
$pr170.Font.Color = 0
$pr170.Font.Italic = $false
$it170 = $c170.Characters($prefix170.Length + 1, $suffix170.Length)
$it170.Font.Name = "Arial"
$it170.Font.Size = 11
$it170.Font.Color = 0
$it170.Font.Italic = $true
$styleSix.Copy()
$c170.PasteSpecial(-4122)

$ws.Cells.Item(170, 1).Value2 = 2022
$ws.Cells.Item(170, 2).Value2 = "May"
$ws.Cells.Item(170, 4).Value2 = "Interview"
$ws.Cells.Item(170, 5).Value2 = "virtual"
$ws.Cells.Item(170, 6).Value2 = "outreach"

# ---- Row 171 ----
$c171 = $ws.Cells.Item(171, 3)
$prefix171 = "UC Santa Barbara SNARL Spring Seminar Series: "
$suffix171 = "Microplastics - a Macro Problem?"
$c171.Value2 = $prefix171 + $suffix171
$pr171 = $c171.Characters(1, $prefix171.Length)
$pr171.Font.Name = "Arial"
$pr171.Font.Size = 11
$pr171.Font.Color = 0
$pr171.Font.Italic = $false
$it171 = $c171.Characters($prefix171.Length + 1, $suffix171.Length)
$it171.Font.Name = "Arial"
$it171.Font.Size = 11
$it171.Font.Color = 0
$it171.Font.Italic = $true
$styleSix.Copy()
$c171.PasteSpecial(-4122)

$ws.Cells.Item(171, 1).Value2 = 2022
$ws.Cells.Item(171, 2).Value2 = "May"
$ws.Cells.Item(171, 4).Value2 = "Platform"
$ws.Cells.Item(171, 5).Value2 = "virtual"
$ws.Cells.Item(171, 6).Value2 = "lecture"

# ---- Row 172 (plain text, no rich runs) ----
$c172 = $ws.Cells.Item(172, 3)
$c172.Value2 = "UCR ENTX Seminar Series: Assessing and Managing Risks of Microplastics"
$styleFive.Copy()
$c172.PasteSpecial(-4122)

$ws.Cells.Item(172, 1).Value2 = 2022
$ws.Cells.Item(172, 2).Value2 = "June"
$ws.Cells.Item(172, 4).Value2 = "Platform"
$ws.Cells.Item(172, 5).Value2 = "in person"
$ws.Cells.Item(172, 6).Value2 = "lecture"

# ---- Row 173 ----
$c173 = $ws.Cells.Item(173, 3)
$prefix173 = "MP Workshop for Early Career Researchers (Athens, Greece):"
$suffix173 = " Assessing and Managing Risks of Microplastics"
$c173.Value2 = $prefix173 + $suffix173
$pr173 = $c173.Characters(1, $prefix173.Length)
$pr173.Font.Name = "Arial"
$pr173.Font.Size = 11
$pr173.Font.Color = 0
$pr173.Font.Italic = $false
$it173 = $c173.Characters($prefix173.Length + 1, $suffix173.Length)
$it173.Font.Name = "Arial"
$it173.Font.Size = 11
$it173.Font.Color = 0
$it173.Font.Italic = $true
$styleSix.Copy()
$c173.PasteSpecial(-4122)

$ws.Cells.Item(173, 1).Value2 = 2022
$ws.Cells.Item(173, 2).Value2 = "June"
$ws.Cells.Item(173, 4).Value2 = "Platform"
$ws.Cells.Item(173, 5).Value2 = "in person"
$ws.Cells.Item(173, 6).Value2 = "lecture"

# Expand the table to cover the new rows
$tbl = $ws.ListObjects.Item(1)
$newRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(173, 6))
$tbl.Resize($newRange)

# Update selection to match the final cursor position after data entry
$ws.Range("D173:F173").Select()

Write-Host "Applied widget rows 168-173"
